# Notes.xlsx -- "ajout de fichier de conversion Excel en xml"
#
# Swap the FirstName/LastName header+data columns (B and C), renumber the
# Note_AP* headers (D1:L1) from the AP2 series (21..29) to the AP1 series
# (11..19), and switch the CNE student-id column (A2:A11) from the old
# 4-digit numbering (2100..2109) to the new 8-digit numbering
# (21000001..21000010). The A3:A11 "+1" formulas are left untouched --
# only the seed value in A2 changes, and the rest recalculate from it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header labels ---
$ws.Range("B1").Value = "FirstName"
$ws.Range("C1").Value = "LastName"

$noteHeaders = @("Note_AP11","Note_AP12","Note_AP13","Note_AP14","Note_AP15","Note_AP16","Note_AP17","Note_AP18","Note_AP19")
$noteCols = @("D","E","F","G","H","I","J","K","L")
for ($i = 0; $i -lt $noteCols.Length; $i++) {
    $ws.Range($noteCols[$i] + "1").Value = $noteHeaders[$i]
}

# --- Row 2: new CNE numbering seed (A3:A11 keep their "+1" formulas) ---
$ws.Range("A2").Value = 21000001

# --- Rows 2-11: FirstName (B) / LastName (C) data swap ---
for ($row = 2; $row -le 11; $row++) {
    $idx = $row - 1
    $ws.Range("B$row").Value = "AP1_FN$idx"
    $ws.Range("C$row").Value = "AP1_LN$idx"
}

# --- Selection moves from I6 to J7 ---
$ws.Range("J7").Select()
